# Apply the changes described in the diff:
#  - Fix typo in the header of column D ("Expexted respins OpenRB" -> "Expexted response OpenRB")
#    (the table column name is driven by the header cell, so it updates automatically)
#  - Fill in a new row (14) describing the "Remove cartridge" command
#  - Update the active selection / top-left visible cell of the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the header cell for column D (this also renames the table column)
$ws.Range("D1").Value = "Expexted response OpenRB"

# Populate the new "Remove cartridge" command row
$ws.Range("A14").Value = "Remove cartridge"
$ws.Range("B14").Value = "REMOVECTRG"
$ws.Range("D14").Value = "CTRG RDY"
$ws.Range("F14").Value = "Lifts all the pinions"

# Update the sheet view: scroll so column B is the left-most visible column,
# and select cell D15 as the active cell
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("D15").Select() | Out-Null
